# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newer scrape snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - row -> new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    4  = 3541
    5  = 3541
    7  = 5065
    8  = 512
    9  = 349
    11 = 688
    13 = 79
    16 = 311
    17 = 33
    21 = 361
    22 = 4892
    23 = 45
    26 = 6016
    29 = 3221
    30 = 330
    31 = 707
    34 = 120
    36 = 1000
    40 = 863
    41 = 967
    42 = 489
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (sheet4.xml) - row -> new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allTypeUpdates = @{
    7  = 3541
    8  = 3541
    10 = 5065
    11 = 512
    12 = 349
    14 = 688
    16 = 79
    19 = 311
    20 = 33
    25 = 361
    26 = 4892
    27 = 45
    30 = 6016
    33 = 3221
    34 = 330
    35 = 707
    39 = 120
    41 = 1000
    45 = 863
    46 = 967
    48 = 491
}
foreach ($row in $allTypeUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allTypeUpdates[$row]
}
